$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.024.21"
$ws.Range("E2").Value = "  -3.36%  "
$ws.Range("D3").Value = "2.337.57"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "85.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.30%  "
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.485"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.80%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.19"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.34%  "
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "2.696.82"
$ws.Range("E13").Value = "  -4.26%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.15%  "
$ws.Range("D16").Value = "2.349.73"
$ws.Range("E16").Value = "  -3.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.759"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "39.998.22"
$ws.Range("E18").Value = "  -3.18%  "
$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.38%  "
$ws.Range("E24").Value = "  -5.57%  "
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("E26").Value = "  -4.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("E28").Value = "  -4.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.91"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0718"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("E36").Value = "  -0.59%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("D42").Value = "1.950.31"
$ws.Range("E42").Value = "  -1.82%  "
$ws.Range("E43").Value = "  -4.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("E45").Value = "  -4.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("E47").Value = "  -5.28%  "
$ws.Range("D48").Value = "2.555.26"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.94%  "
